$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.114.21"
$ws.Range("E2").Value = "  +3.36%  "

$ws.Range("D3").Value = "1.597.71"
$ws.Range("E3").Value = "  +1.99%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.15"
$ws.Range("E5").Value = "  +2.20%  "

$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("E7").Value = "  +1.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  +2.37%  "

$ws.Range("E9").Value = "  +1.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.89"
$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("E11").Value = "  +4.99%  "

$ws.Range("D12").Value = "1.822.56"
$ws.Range("E12").Value = "  +2.16%  "

$ws.Range("D13").Value = "1.601.50"
$ws.Range("E13").Value = "  +2.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").Value = "26.085.23"
$ws.Range("E16").Value = "  +3.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.42"
$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("D18").Value = "0.0₃0720"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "204.20"
$ws.Range("E20").Value = "  +9.90%  "

$ws.Range("E21").Value = "  +2.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.31"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("E23").Value = "  +1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.82"
$ws.Range("E24").Value = "  +11.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.87"
$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("E27").Value = "  -2.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.23"
$ws.Range("E28").Value = "  +2.47%  "

$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("E32").Value = "  +2.65%  "

$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("E34").Value = "  +0.75%  "

$ws.Range("E35").Value = "  +2.68%  "

$ws.Range("E36").Value = "  +9.39%  "

$ws.Range("D37").Value = "1.105.39"
$ws.Range("E37").Value = "  +1.79%  "

$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("E40").Value = "  +0.44%  "

$ws.Range("E41").Value = "  -0.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.777"
$ws.Range("E42").Value = "  +1.74%  "

$ws.Range("D43").Value = "1.736.09"
$ws.Range("E43").Value = "  +2.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "92.62"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("E45").Value = "  +0.53%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.51"
$ws.Range("E46").Value = "  +6.75%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0103"
$ws.Range("E47").Value = "  -4.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.31"
$ws.Range("E48").Value = "  +1.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0505"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("E50").Value = "  +0.77%  "

$ws.Range("E51").Value = "  +0.09%  "

